$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text in B1: "Módulo [A]" -> "Modulo [A]" (accent removed)
$ws.Range("B1").Value = "Modulo [A]"

# Zero out the "Modulo [A]" column values from row 4 through row 52
for ($r = 4; $r -le 52; $r++) {
    $ws.Cells.Item($r, 2).Value = 0
}

# Update the active selection to B9:B52 with active cell B9
$ws.Range("B9:B52").Select()
